$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text formatting
# (values like "1.00", "0.0691", "375.00" would otherwise be
# auto-converted to numbers by Excel, losing trailing zeros / exact text)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.517.71"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "3.440.07"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "580.82"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("D6").Value = "175.12"
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("D9").Value = "3.434.94"
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").Value = "6.82"
$ws.Range("E11").Value = "  -3.65%  "
$ws.Range("D12").Value = "0.418"
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("D13").Value = "4.036.93"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "30.95"
$ws.Range("E14").Value = "  -3.68%  "
$ws.Range("D15").Value = "0.132"
$ws.Range("E15").Value = "  -3.42%  "
$ws.Range("D16").Value = "66.557.62"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "0.0000171"
$ws.Range("E17").Value = "  -2.95%  "
$ws.Range("D18").Value = "3.445.58"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").Value = "6.00"
$ws.Range("E19").Value = "  -3.96%  "
$ws.Range("D20").Value = "13.77"
$ws.Range("E20").Value = "  -3.58%  "
$ws.Range("D21").Value = "375.00"
$ws.Range("E21").Value = "  -3.88%  "
$ws.Range("D22").Value = "7.67"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").Value = "0.995"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "5.72"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "70.75"
$ws.Range("E25").Value = "  -2.84%  "
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("E27").Value = "  -2.32%  "
$ws.Range("D28").Value = "9.86"
$ws.Range("E28").Value = "  -4.70%  "
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -5.20%  "
$ws.Range("D32").Value = "23.78"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("E33").Value = "  -3.00%  "
$ws.Range("E34").Value = "  -6.27%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "7.02"
$ws.Range("E36").Value = "  -4.53%  "
$ws.Range("E37").Value = "  -5.71%  "
$ws.Range("D38").Value = "158.93"
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("D39").Value = "0.874"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "27.01"
$ws.Range("E40").Value = "  +4.08%  "
$ws.Range("E41").Value = "  -4.78%  "
$ws.Range("E42").Value = "  -4.20%  "
$ws.Range("D43").Value = "6.50"
$ws.Range("E43").Value = "  -5.18%  "
$ws.Range("E44").Value = "  -4.01%  "
$ws.Range("D45").Value = "2.685.12"
$ws.Range("E45").Value = "  -5.56%  "
$ws.Range("D46").Value = "0.0691"
$ws.Range("E46").Value = "  -4.45%  "
$ws.Range("D47").Value = "25.16"
$ws.Range("E47").Value = "  -4.99%  "
$ws.Range("D48").Value = "40.34"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("D50").Value = "318.80"
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -4.40%  "
